{"js": "// Load all paragraphs in the main body so we can locate the anchor points by\n// their current text content (robust to the exact paragraph count / rsids\n// already present in the document).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Small helper: find the (first) paragraph item whose text equals `text`,\n// starting the scan at `startAt` (defaults to 0). Throws if not found so a\n// mismatch fails loudly instead of silently doing nothing.\nfunction findParagraph(text, startAt) {\n  const items = paragraphs.items;\n  for (let i = startAt || 0; i < items.length; i++) {\n    if (items[i].text === text) return items[i];\n  }\n  throw new Error(\"Paragraph not found: \" + text);\n}\n\n// --- Change 1 -------------------------------------------------------------\n// Insert \"# Initialiser un compteur de lignes\" / \"line_number = 1\" / a blank\n// paragraph just before the \"# Parcourir chaque ligne ...\" comment.\nconst parcourirPara = findParagraph(\n  \"# Parcourir chaque ligne dans la colonne A (A2 \u00e0 A... en supposant que A1 est un en-t\u00eate)\"\n);\nconst initCounterPara = parcourirPara.insertParagraph(\n  \"# Initialiser un compteur de lignes\",\n  \"Before\"\n);\nconst lineNumberInitPara = initCounterPara.insertParagraph(\"line_number = 1\", \"After\");\nlineNumberInitPara.insertParagraph(\"\", \"After\");\n\n// --- Change 2 -------------------------------------------------------------\n// Drop the `values_only=True` kwarg from the iter_rows(...) call.\nconst forRowPara = findParagraph(\n  \"for row in sheet.iter_rows(min_row=2, min_col=1, max_col=1, values_only=True):\"\n);\nforRowPara.insertText(\n  \"for row in sheet.iter_rows(min_row=2, min_col=1, max_col=1):\",\n  \"Replace\"\n);\n\n// --- Change 3 -------------------------------------------------------------\n// Split \"xml_content = row[0]  # ...\" into \"cell = row[0]\" followed by\n// \"xml_content = cell.value  # ...\".\nconst xmlContentPara = findParagraph(\n  \"    xml_content = row[0]  # R\u00e9cup\u00e9rer la valeur de la cellule en colonne A\"\n);\nxmlContentPara.insertText(\"    cell = row[0]\", \"Replace\");\nxmlContentPara.insertParagraph(\n  \"    xml_content = cell.value  # R\u00e9cup\u00e9rer la valeur de la cellule en colonne A\",\n  \"After\"\n);\n\n// --- Change 4 -------------------------------------------------------------\n// Remove the \"# Obtenir l'index de la ligne\" / \"row_num = sheet._current_row\"\n// / blank-line block, and update the two paragraphs that follow it.\nconst obtenirIndexPara = findParagraph(\"        # Obtenir l'index de la ligne\");\nconst rowNumPara = findParagraph(\"        row_num = sheet._current_row\");\nconst blankAfterRowNum = findParagraph(\"\", 0); // placeholder, replaced below\n\n// Locate the specific blank paragraph that sits between \"row_num = ...\" and\n// the \"# Nom du fichier ...\" comment by scanning paragraphs around rowNumPara.\n{\n  const items = paragraphs.items;\n  const rowNumIdx = items.indexOf(rowNumPara);\n  const blankPara = items[rowNumIdx + 1];\n  obtenirIndexPara.delete();\n  rowNumPara.delete();\n  blankPara.delete();\n}\n\nconst filenameCommentPara = findParagraph(\n  '        # Nom du fichier .xml (par exemple : ligne 2 devient \"ligne_2.xml\")'\n);\nfilenameCommentPara.insertText(\n  \"        # Nom du fichier .xml bas\u00e9 sur le num\u00e9ro de ligne\",\n  \"Replace\"\n);\n\nconst xmlFilenamePara = findParagraph('        xml_filename = f\"ligne_{row_num}.xml\"');\nxmlFilenamePara.insertText(\n  '        xml_filename = f\"ligne_{line_number}.xml\"',\n  \"Replace\"\n);\n\n// --- Change 5 -------------------------------------------------------------\n// After the \"print(f\\\"Fichier cr\u00e9\u00e9 : {output_path}\\\")\" line, add a blank\n// paragraph, a comment, and the counter increment.\nconst fichierCreePara = findParagraph(\n  '        print(f\"Fichier cr\u00e9\u00e9 : {output_path}\")'\n);\nconst blankAfterPrint = fichierCreePara.insertParagraph(\"\", \"After\");\nconst incrementCommentPara = blankAfterPrint.insertParagraph(\n  \"    # Incr\u00e9menter le compteur de lignes\",\n  \"After\"\n);\nincrementCommentPara.insertParagraph(\"    line_number += 1\", \"After\");\n\nawait context.sync();\n", "ps1": "# Helper: return the 1-based index of the first paragraph in $doc whose\n# text (paragraph mark stripped) exactly equals $text. Throws if missing.\nfunction Find-ParaIndex {\n    param($doc, $text)\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $t = $doc.Paragraphs.Item($i).Range.Text\n        $t2 = $t.TrimEnd([char]13, [char]7)\n        if ($t2 -eq $text) {\n            return $i\n        }\n    }\n    throw \"Paragraph not found: $text\"\n}\n\n$d = $word.ActiveDocument\n\n# --- Change 1 ---------------------------------------------------------\n# Insert \"# Initialiser un compteur de lignes\" / \"line_number = 1\" / a\n# blank paragraph just before the \"# Parcourir chaque ligne ...\" comment.\n$idx = Find-ParaIndex $d \"# Parcourir chaque ligne dans la colonne A (A2 \u00e0 A... en supposant que A1 est un en-t\u00eate)\"\n$d.Paragraphs.Item($idx).Range.InsertParagraphBefore()\n$d.Paragraphs.Item($idx).Range.Text = \"# Initialiser un compteur de lignes\"\n$idxLineNumber = $idx + 1\n$d.Paragraphs.Item($idxLineNumber).Range.InsertParagraphBefore()\n$d.Paragraphs.Item($idxLineNumber).Range.Text = \"line_number = 1\"\n$idxBlank = $idxLineNumber + 1\n$d.Paragraphs.Item($idxBlank).Range.InsertParagraphBefore()\n# (the paragraph left at $idxBlank-1 is the new blank line; nothing else to set)\n\n# --- Change 2 ---------------------------------------------------------\n# Drop the `values_only=True` kwarg from the iter_rows(...) call.\n$idxForRow = Find-ParaIndex $d \"for row in sheet.iter_rows(min_row=2, min_col=1, max_col=1, values_only=True):\"\n$d.Paragraphs.Item($idxForRow).Range.Text = \"for row in sheet.iter_rows(min_row=2, min_col=1, max_col=1):\"\n\n# --- Change 3 ---------------------------------------------------------\n# Split \"xml_content = row[0]  # ...\" into \"cell = row[0]\" followed by a\n# new \"xml_content = cell.value  # ...\" paragraph.\n$idxXmlContent = Find-ParaIndex $d \"    xml_content = row[0]  # R\u00e9cup\u00e9rer la valeur de la cellule en colonne A\"\n$d.Paragraphs.Item($idxXmlContent).Range.Text = \"    cell = row[0]\"\n$d.Paragraphs.Item($idxXmlContent).Range.InsertParagraphAfter()\n$d.Paragraphs.Item($idxXmlContent + 1).Range.Text = \"    xml_content = cell.value  # R\u00e9cup\u00e9rer la valeur de la cellule en colonne A\"\n\n# --- Change 4 ---------------------------------------------------------\n# Remove the \"# Obtenir l'index de la ligne\" / \"row_num = sheet._current_row\"\n# / blank-line block (3 consecutive paragraphs), then update the two\n# paragraphs that follow it.\n$idxObtenir = Find-ParaIndex $d \"        # Obtenir l'index de la ligne\"\n$d.Paragraphs.Item($idxObtenir).Range.Delete()\n$d.Paragraphs.Item($idxObtenir).Range.Delete()\n$d.Paragraphs.Item($idxObtenir).Range.Delete()\n\n$idxFilenameComment = Find-ParaIndex $d \"        # Nom du fichier .xml (par exemple : ligne 2 devient `\"ligne_2.xml`\")\"\n$d.Paragraphs.Item($idxFilenameComment).Range.Text = \"        # Nom du fichier .xml bas\u00e9 sur le num\u00e9ro de ligne\"\n\n$idxFilename = Find-ParaIndex $d \"        xml_filename = f`\"ligne_{row_num}.xml`\"\"\n$d.Paragraphs.Item($idxFilename).Range.Text = \"        xml_filename = f`\"ligne_{line_number}.xml`\"\"\n\n# --- Change 5 ---------------------------------------------------------\n# After the \"print(f\\\"Fichier cr\u00e9\u00e9 : {output_path}\\\")\" line, add a blank\n# paragraph, a comment, and the counter increment.\n$idxPrint = Find-ParaIndex $d \"        print(f`\"Fichier cr\u00e9\u00e9 : {output_path}`\")\"\n$d.Paragraphs.Item($idxPrint).Range.InsertParagraphAfter()\n$idxNewBlank = $idxPrint + 1\n$d.Paragraphs.Item($idxNewBlank).Range.InsertParagraphAfter()\n$idxIncComment = $idxNewBlank + 1\n$d.Paragraphs.Item($idxIncComment).Range.Text = \"    # Incr\u00e9menter le compteur de lignes\"\n$d.Paragraphs.Item($idxIncComment).Range.InsertParagraphAfter()\n$idxIncrement = $idxIncComment + 1\n$d.Paragraphs.Item($idxIncrement).Range.Text = \"    line_number += 1\"\n"}
